$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The reference importer now logs (instead of silently skipping)
# rows whose "Doi" field isn't a real "10.xxx" DOI, and the sample
# data was re-exported: the record that used to sit at row 2 (its
# "Doi" column actually holds a title, not a DOI - "Retrospective
# Clinical Analysis...") moves down to row 3, and the record that
# used to sit at row 3 (a real "10.xxx" DOI, with its hyperlink and
# parsed metadata) moves up to become the new row 2.
#
# Swap the two rows' content *and* formatting using Copy/Paste
# through scratch cells well outside the used range, so each
# destination keeps the exact style the data should carry.
# -----------------------------------------------------------------

# 1) Stash full old row 2 (A:H) and old row 3 (A:B) in scratch cells.
$ws.Range("A2:H2").Copy($ws.Range("A200:H200"))
$ws.Range("A3:B3").Copy($ws.Range("A201:B201"))

# Also stash old row 3's A-cell format alone (the underlined /
# hyperlink look) so it can be re-applied after the hyperlink API
# clobbers it with its own built-in "Hyperlink" style below.
$ws.Range("A3").Copy()
$ws.Range("A202").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 2) Clear the two source rows so nothing lingers in unused columns.
$ws.Range("A2:H2").Clear()
$ws.Range("A3:H3").Clear()

# 3) Old row 3 (A:B) becomes the new row 2 (A:B) - DOI + date, with
#    the hyperlink-style formatting and the mmmm-yyyy date format.
$ws.Range("A201:B201").Copy($ws.Range("A2:B2"))

# 4) Old row 2's A:B (title + year) becomes the new row 3 (A:B).
$ws.Range("A200:B200").Copy($ws.Range("A3:B3"))

# 5) Old row 2's C:H (the "Unknown Title" / "not found" / etc block)
#    stays attached to the DOI record, so it comes back on row 2.
$ws.Range("C200:H200").Copy($ws.Range("C2:H2"))

# 6) Clean up the scratch area used for the row swap.
$ws.Range("A200:H201").Clear()

# -----------------------------------------------------------------
# Move the hyperlink from the old A3 to the new A2.
# -----------------------------------------------------------------
$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "about:blank") | Out-Null

# Adding a hyperlink registers Excel's built-in "Hyperlink" cell
# style; drop it again and restore the original underline format
# so A2 keeps looking exactly like the old (hyperlinked) A3 did.
$wb.Styles.Item("Hyperlink").Delete()
$ws.Range("A202").Copy()
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A202").Clear()

# -----------------------------------------------------------------
# A couple of stray formatted-but-empty cells also show up in the
# edited sheet (D6 and B13), carrying the same fill/format as the
# other "year" cells in column B. Paste just the formatting over.
# -----------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("D6").PasteSpecial(-4122) | Out-Null

$ws.Range("B4").Copy()
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
